$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B28: was stored as text "4", should become a real number 4.
$ws.Range("B28").Value = 4

# Insert the new row 29 with the split-out annotation row.
$ws.Range("A29").Value = "Sunsi Wu"

# B29 must stay a text "4" (matching the source data's original quirk),
# not get auto-coerced into a number by Excel's type inference.
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "4"
$ws.Range("B29").ClearFormats()

$ws.Range("C29").Value = "elaborate"
$ws.Range("D29").Value = "ACK"
$ws.Range("E29").Value = "OTH"
$ws.Range("F29").Value = "9b81a0cf-ae6f-4476-b619-1b75e1becf94"
$ws.Range("G29").Value = "B1ae1lZRb_annotated.xlsx"
$ws.Range("H29").Value = "We will elaborate on this aspect in the final version of the paper."
